# Apply the "moving average filtering" progress update + journal entry.
$wb = $excel.ActiveWorkbook

# --- Progress sheet: mark the "Max Average Filtering" task row (row 12) done ---
$wsProgress = $wb.Worksheets.Item("Progress")
$wsProgress.Range("G12").Value = 45734
# Copy G12's (date) formatting onto H12 before writing the status text, so the
# status cell picks up the same border/number-format styling the author left it with.
$wsProgress.Range("G12").Copy()
$wsProgress.Range("H12").PasteSpecial(-4122)
$wsProgress.Range("H12").Value = "done"

# --- Journal sheet: add a new dated entry about the moving average filtering work ---
$wsJournal = $wb.Worksheets.Item("Journal")
$wsJournal.Range("A8").Value = "Moving Average Filtering"
$wsJournal.Range("B8").Value = 45734
$wsJournal.Range("C8").Value = "Incountered unmatching result while creating this, where the reason was due to skiping of 1st index value during the calculation. This is connected to The adding of 1 in front of the output signals in NEO transform. If NEO transform later gets rid of the 1 value in the front of the signal, this may need to be fixed with it.(moving_average_filtering function in activation_detection.c)"
$wsJournal.Rows.Item(8).RowHeight = 72.9

# Update the on-screen selections to match where the author ended up, Journal
# first so the final active sheet/selection ends up back on Progress.
$wsJournal.Range("D7").Select()
$wsProgress.Range("L12").Select()
